$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.473.76'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.572.49'
$ws.Range("E3").Value = '  +0.32%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.73'
$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3703'
$ws.Range("E7").Value = '  -1.85%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.93'
$ws.Range("E8").Value = '  +1.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3382'
$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.143'
$ws.Range("E10").Value = '  +0.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07544'
$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.23'
$ws.Range("E13").Value = '  +0.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.026'
$ws.Range("E14").Value = '  +0.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.959'
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.591.35'
$ws.Range("E16").Value = '  +1.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001119'
$ws.Range("E17").Value = '  -1.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.73'
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.293'
$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.41'
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.19'
$ws.Range("E23").Value = '  +1.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.493.36'
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.369'
$ws.Range("E25").Value = '  -1.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.600'
$ws.Range("E26").Value = '  -3.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.05'
$ws.Range("E27").Value = '  -0.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.16'
$ws.Range("E28").Value = '  +1.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.061'
$ws.Range("E29").Value = '  +0.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.20'
$ws.Range("E30").Value = '  -0.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.764.67'
$ws.Range("E31").Value = '  +1.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.076'
$ws.Range("E32").Value = '  +8.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.245'
$ws.Range("E33").Value = '  +2.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.010'
$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.777'
$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08354'
$ws.Range("E36").Value = '  -1.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02483'
$ws.Range("E37").Value = '  -1.28%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.355'
$ws.Range("E38").Value = '  -4.64%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2301'
$ws.Range("E39").Value = '  +0.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06535'
$ws.Range("E40").Value = '  +0.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.443'
$ws.Range("E41").Value = '  +0.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.36'
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6216'
$ws.Range("E43").Value = '  -1.67%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.99'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.808'
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5853'
$ws.Range("E47").Value = '  -1.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.28'
$ws.Range("E48").Value = '  +3.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.069'
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("E50").Value = '  -2.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07332'
$ws.Range("E51").Value = '  +0.12%  '
